$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '285.71'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.83%'
$ws.Range('G2').NumberFormat = '@'
$ws.Range('G2').Value = '13'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '29.46'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '4.13%'
$ws.Range('G3').NumberFormat = '@'
$ws.Range('G3').Value = '13'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.080'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.22%'
$ws.Range('G4').NumberFormat = '@'
$ws.Range('G4').Value = '13'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.06708'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '2.67%'
$ws.Range('G5').NumberFormat = '@'
$ws.Range('G5').Value = '13'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '7.323'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '1.13%'
$ws.Range('G6').NumberFormat = '@'
$ws.Range('G6').Value = '13'
$ws.Range('B7').NumberFormat = '@'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').NumberFormat = '@'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.446'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '2.22%'
$ws.Range('G7').NumberFormat = '@'
$ws.Range('G7').Value = '13'
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.371'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-8.68%'
$ws.Range('G8').NumberFormat = '@'
$ws.Range('G8').Value = '13'
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9030'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-1.73%'
$ws.Range('G9').NumberFormat = '@'
$ws.Range('G9').Value = '13'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1577'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '2.09%'
$ws.Range('G10').NumberFormat = '@'
$ws.Range('G10').Value = '13'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07088'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '9.75%'
$ws.Range('G11').NumberFormat = '@'
$ws.Range('G11').Value = '13'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07603'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '0.63%'
$ws.Range('G12').NumberFormat = '@'
$ws.Range('G12').Value = '13'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.02924'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '6.01%'
$ws.Range('G13').NumberFormat = '@'
$ws.Range('G13').Value = '13'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.08993'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '0.30%'
$ws.Range('G14').NumberFormat = '@'
$ws.Range('G14').Value = '13'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001589'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.12%'
$ws.Range('G15').NumberFormat = '@'
$ws.Range('G15').Value = '13'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.04500'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '2.02%'
$ws.Range('G16').NumberFormat = '@'
$ws.Range('G16').Value = '13'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'One'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0006490'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '1.09%'
$ws.Range('G17').NumberFormat = '@'
$ws.Range('G17').Value = '13'
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'TigerCash'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.006202'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '0.62%'
$ws.Range('G18').NumberFormat = '@'
$ws.Range('G18').Value = '13'
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'LEO'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.451'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '0.10%'
$ws.Range('G19').NumberFormat = '@'
$ws.Range('G19').Value = '13'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-0.41%'
$ws.Range('G20').NumberFormat = '@'
$ws.Range('G20').Value = '13'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.3232'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '1.14%'
$ws.Range('G21').NumberFormat = '@'
$ws.Range('G21').Value = '13'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '2.05%'
$ws.Range('G22').NumberFormat = '@'
$ws.Range('G22').Value = '13'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.952'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-0.32%'
$ws.Range('G23').NumberFormat = '@'
$ws.Range('G23').Value = '13'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '3.58%'
$ws.Range('G24').NumberFormat = '@'
$ws.Range('G24').Value = '13'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.001206'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '1.87%'
$ws.Range('G25').NumberFormat = '@'
$ws.Range('G25').Value = '13'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-1.40%'
$ws.Range('G26').NumberFormat = '@'
$ws.Range('G26').Value = '13'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-3.79%'
$ws.Range('G27').NumberFormat = '@'
$ws.Range('G27').Value = '13'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0001621'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '0.33%'
$ws.Range('G28').NumberFormat = '@'
$ws.Range('G28').Value = '13'
$ws.Range('G29').NumberFormat = '@'
$ws.Range('G29').Value = '13'
$ws.Range('G30').NumberFormat = '@'
$ws.Range('G30').Value = '13'
$ws.Range('G31').NumberFormat = '@'
$ws.Range('G31').Value = '13'
$ws.Range('G32').NumberFormat = '@'
$ws.Range('G32').Value = '13'
$ws.Range('G33').NumberFormat = '@'
$ws.Range('G33').Value = '13'
$ws.Range('G34').NumberFormat = '@'
$ws.Range('G34').Value = '13'
$ws.Range('G35').NumberFormat = '@'
$ws.Range('G35').Value = '13'
$ws.Range('G36').NumberFormat = '@'
$ws.Range('G36').Value = '13'
$ws.Range('G37').NumberFormat = '@'
$ws.Range('G37').Value = '13'
$ws.Range('G38').NumberFormat = '@'
$ws.Range('G38').Value = '13'
$ws.Range('G39').NumberFormat = '@'
$ws.Range('G39').Value = '13'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04234'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '2.20%'
$ws.Range('G40').NumberFormat = '@'
$ws.Range('G40').Value = '13'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006794'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '1.85%'
$ws.Range('G41').NumberFormat = '@'
$ws.Range('G41').Value = '13'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1239'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '0.73%'
$ws.Range('G42').NumberFormat = '@'
$ws.Range('G42').Value = '13'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002232'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '4.55%'
$ws.Range('G43').NumberFormat = '@'
$ws.Range('G43').Value = '13'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '4.40%'
$ws.Range('G44').NumberFormat = '@'
$ws.Range('G44').Value = '13'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005783'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '4.17%'
$ws.Range('G45').NumberFormat = '@'
$ws.Range('G45').Value = '13'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.967'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.06%'
$ws.Range('G46').NumberFormat = '@'
$ws.Range('G46').Value = '13'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '15.68%'
$ws.Range('G47').NumberFormat = '@'
$ws.Range('G47').Value = '13'
$ws.Range('G48').NumberFormat = '@'
$ws.Range('G48').Value = '13'
$ws.Range('G49').NumberFormat = '@'
$ws.Range('G49').Value = '13'
$ws.Range('G50').NumberFormat = '@'
$ws.Range('G50').Value = '13'
$ws.Range('G51').NumberFormat = '@'
$ws.Range('G51').Value = '13'
